$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'230.42"

$ws.Range("D3").Value = "'22.67"

$ws.Range("D4").Value = "'5.294"

$ws.Range("D6").Value = "'3.381"

$ws.Range("D8").Value = "'1.063"

$ws.Range("D9").Value = "'0.7825"

$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").Value = "'0.01168"
$ws.Range("E10").Value = "9OneONEBestin24h"

$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1399"
$ws.Range("E11").Value = "10WazirXWRX"

$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.07388"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"

$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D13").Value = "'0.03176"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"

$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "'0.02974"
$ws.Range("E14").Value = "13BitrueCoinBTR"

$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "'0.09258"
$ws.Range("E15").Value = "14BitMartTokenBMX"

$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001680"
$ws.Range("E16").Value = "15BitForexTokenBF"

$ws.Range("B17").Value = "MCDex"
$ws.Range("C17").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D17").Value = "'3.257"
$ws.Range("E17").Value = "16MCDexMCB"

$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D18").Value = "'0.04742"
$ws.Range("E18").Value = "17CoinExTokenCET"

$ws.Range("D19").Value = "'0.006241"

$ws.Range("D20").Value = "'0.005221"

$ws.Range("D22").Value = "'0.0001501"

$ws.Range("D26").Value = "'0.1253"

$ws.Range("D27").Value = "'0.0004994"
$ws.Range("E27").Value = "26UpBotsUBXT"

$ws.Range("D40").Value = "'0.04035"

$ws.Range("D41").Value = "'0.007009"

$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.003503"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").Value = "'0.1039"
$ws.Range("E43").Value = "42BKEXTokenBKK"

$ws.Range("D44").Value = "'0.01034"

$ws.Range("D45").Value = "'0.00005444"

$ws.Range("D46").Value = "'0.00000000751"

$ws.Range("D47").Value = "'0.7858"

$ws.Range("D48").Value = "'0.04074"

$ws.Range("D49").Value = "'0.00002102"
